{"js": "// Simplified Chinese -> Traditional Chinese (zh-TW) translation update\n// for \"Partner email \u2013 document verification failed\" template, plus a\n// couple of English runs that were translated into Traditional Chinese.\n//\n// Each pair is [searchText, replacementText]. We search the whole body\n// for (mostly unique) runs of text and replace them in place, which\n// preserves the surrounding run formatting (bold/color/highlight/etc.).\n\nconst replacements = [\n  [\"\u82f1\u8bed\", \"\u82f1\u8a9e\"],\n  [\" / \u8461\u8404\u7259\u8bed / \u6cd5\u8bed / \u6cf0\u8bed / \u8d8a\u5357\u8bed / \u897f\u73ed\u7259\u8bed\", \" / \u8461\u8404\u7259\u8a9e / \u6cd5\u8a9e / \u6cf0\u8a9e / \u8d8a\u5357\u8a9e / \u897f\u73ed\u7259\u8a9e\"],\n  [\"\u7b80\u4ecb\", \"\u7c21\u4ecb\"],\n  [\n    \"\u4e00\u5c01\u53d1\u9001\u7ed9\u76ee\u6807\u56fd\u5bb6\u4e2d\u672a\u901a\u8fc7\u6211\u4eec\u9a8c\u8bc1\u6d41\u7a0b\u7684\u5408\u4f5c\u4f19\u4f34\u7684\u7535\u5b50\u90ae\u4ef6\u3002 \u5c06\u901a\u8fc7 customer.io \u53d1\u9001\",\n    \"\u767c\u9001\u7d66\u76ee\u6a19\u570b\u5bb6\u4e2d\u90a3\u4e9b\u6587\u4ef6\u672a\u901a\u904e\u6211\u5011\u9a57\u8b49\u6d41\u7a0b\u7684\u5408\u4f5c\u5925\u4f34\u7684\u96fb\u5b50\u90f5\u4ef6\u3002 \u5c07\u901a\u904e customer.io \u767c\u9001\",\n  ],\n  [\"\u76ee\u6807\u53d7\u4f17\", \"\u76ee\u6a19\u53d7\u773e\"],\n  [\"\u63d0\u4ea4\u4e86\u9519\u8bef/\u4e0d\u5b8c\u6574\u6587\u4ef6\u7684\u9080\u8bf7\u5408\u4f5c\u4f19\u4f34\", \"\u63d0\u4ea4\u932f\u8aa4/\u4e0d\u5b8c\u6574\u6587\u6a94\u7684\u88ab\u9080\u8acb\u5408\u4f5c\u5925\u4f34\"],\n  [\"\u4e3b\u9898\u884c\", \"\u4e3b\u984c\u884c\"],\n  [\"[\u4e8b\u4ef6\u540d\u79f0]\", \"[\u4e8b\u4ef6\u540d\u7a31]\"],\n  [\" \u2014 \u6587\u6863\u9a8c\u8bc1\u5931\u8d25 \", \" \u2014 \u6587\u4ef6\u9a57\u8b49\u5931\u6557 \"],\n  [\"\u554a\u54e6\uff01 \u6587\u4ef6\u65e0\u6cd5\u9a8c\u8bc1\", \"\u554a\u54e6\uff01 \u6587\u6a94\u7121\u6cd5\u9a57\u8b49\"],\n  [\"[\u5408\u4f5c\u4f19\u4f34\u59d3\u540d]\", \"[\u5408\u4f5c\u5925\u4f34\u59d3\u540d]\"],\n  [\n    \"We regret to inform you that your documents have failed our verification process as we found the following issues with them: \",\n    \"\u5f88\u907a\u61be\u5730\u901a\u77e5\u60a8\uff0c\u60a8\u7684\u6587\u6a94\u672a\u901a\u904e\u9a57\u8b49\u6d41\u7a0b\uff0c\u56e0\u70ba\u6211\u5011\u767c\u73fe\u4ee5\u4e0b\u554f\u984c\uff1a \",\n  ],\n  [\"\u75ab\u82d7\u63a5\u79cd\u8bc1\u4e66\u526f\u672c\", \"\u60a8\u7684\u75ab\u82d7\u63a5\u7a2e\u8b49\u660e\u526f\u672c\"],\n  [\": \u6587\u4ef6\u4e0d\u6e05\u695a\", \": \u6587\u6a94\u4e0d\u6e05\u695a\"],\n  [\"[\u6587\u4ef6 2]\", \"[\u6587\u6a94 2]\"],\n  [\": [problem]\", \": [\u554f\u984c]\"],\n  [\"\u8bf7\u5728 \", \"\u8acb\u5728 \"],\n  [\" \u4e4b\u524d\u91cd\u65b0\u63d0\u4ea4\u4e0a\u8ff0\u6587\u4ef6\uff0c\u4ee5\u4fbf\u6211\u4eec\u8fdb\u884c\u5fc5\u8981\u7684\u5b89\u6392\u3002\", \" \u4e4b\u524d\u91cd\u65b0\u63d0\u4ea4\u4e0a\u8ff0\u6587\u6a94\uff0c\u4ee5\u4fbf\u6211\u5011\u9032\u884c\u5fc5\u8981\u7684\u5b89\u6392\u3002\"],\n  [\"\u5982\u6709\u4efb\u4f55\u7591\u95ee\uff0c\u8bf7\u901a\u8fc7 \", \"\u5982\u6709\u4efb\u4f55\u7591\u554f\uff0c\u8acb\u901a\u904e \"],\n  [\"[\u7535\u5b50\u90ae\u4ef6\u5730\u5740]\", \"[\u96fb\u5b50\u90f5\u4ef6\u5730\u5740]\"],\n  [\"[WHATSAPP \u53f7\u7801]\", \"[WHATSAPP \u865f\u78bc]\"],\n  [\" (WhatsApp) \u8054\u7cfb\u60a8\u7684\u533a\u57df\u7ecf\u7406 \", \" (WhatsApp) \u806f\u7e6b\u60a8\u7684\u5340\u57df\u7d93\u7406, \"],\n  [\"[NAME]\", \"[\u59d3\u540d]\"],\n  [\" \u3002 \", \"\u3002 \"],\n];\n\nfor (const [searchText, replacementText] of replacements) {\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacementText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Simplified Chinese -> Traditional Chinese (zh-TW) translation update\n# for \"Partner email \u2013 document verification failed\" template, plus a\n# couple of English runs that were translated into Traditional Chinese.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @(\"\u82f1\u8bed\", \"\u82f1\u8a9e\"),\n    # Leading space is unchanged by the edit; starting the match exactly at\n    # the hyperlink/run boundary makes Word's Find/Replace inherit the\n    # adjacent hyperlink's formatting, so the (unchanged) space is left out.\n    @(\"\u8461\u8404\u7259\u8bed / \u6cd5\u8bed / \u6cf0\u8bed / \u8d8a\u5357\u8bed / \u897f\u73ed\u7259\u8bed\", \"\u8461\u8404\u7259\u8a9e / \u6cd5\u8a9e / \u6cf0\u8a9e / \u8d8a\u5357\u8a9e / \u897f\u73ed\u7259\u8a9e\"),\n    @(\"\u7b80\u4ecb\", \"\u7c21\u4ecb\"),\n    @(\"\u4e00\u5c01\u53d1\u9001\u7ed9\u76ee\u6807\u56fd\u5bb6\u4e2d\u672a\u901a\u8fc7\u6211\u4eec\u9a8c\u8bc1\u6d41\u7a0b\u7684\u5408\u4f5c\u4f19\u4f34\u7684\u7535\u5b50\u90ae\u4ef6\u3002 \u5c06\u901a\u8fc7 customer.io \u53d1\u9001\", \"\u767c\u9001\u7d66\u76ee\u6a19\u570b\u5bb6\u4e2d\u90a3\u4e9b\u6587\u4ef6\u672a\u901a\u904e\u6211\u5011\u9a57\u8b49\u6d41\u7a0b\u7684\u5408\u4f5c\u5925\u4f34\u7684\u96fb\u5b50\u90f5\u4ef6\u3002 \u5c07\u901a\u904e customer.io \u767c\u9001\"),\n    @(\"\u76ee\u6807\u53d7\u4f17\", \"\u76ee\u6a19\u53d7\u773e\"),\n    @(\"\u63d0\u4ea4\u4e86\u9519\u8bef/\u4e0d\u5b8c\u6574\u6587\u4ef6\u7684\u9080\u8bf7\u5408\u4f5c\u4f19\u4f34\", \"\u63d0\u4ea4\u932f\u8aa4/\u4e0d\u5b8c\u6574\u6587\u6a94\u7684\u88ab\u9080\u8acb\u5408\u4f5c\u5925\u4f34\"),\n    @(\"\u4e3b\u9898\u884c\", \"\u4e3b\u984c\u884c\"),\n    @(\"[\u4e8b\u4ef6\u540d\u79f0]\", \"[\u4e8b\u4ef6\u540d\u7a31]\"),\n    @(\" \u2014 \u6587\u6863\u9a8c\u8bc1\u5931\u8d25 \", \" \u2014 \u6587\u4ef6\u9a57\u8b49\u5931\u6557 \"),\n    @(\"\u554a\u54e6\uff01 \u6587\u4ef6\u65e0\u6cd5\u9a8c\u8bc1\", \"\u554a\u54e6\uff01 \u6587\u6a94\u7121\u6cd5\u9a57\u8b49\"),\n    @(\"[\u5408\u4f5c\u4f19\u4f34\u59d3\u540d]\", \"[\u5408\u4f5c\u5925\u4f34\u59d3\u540d]\"),\n    @(\"We regret to inform you that your documents have failed our verification process as we found the following issues with them: \", \"\u5f88\u907a\u61be\u5730\u901a\u77e5\u60a8\uff0c\u60a8\u7684\u6587\u6a94\u672a\u901a\u904e\u9a57\u8b49\u6d41\u7a0b\uff0c\u56e0\u70ba\u6211\u5011\u767c\u73fe\u4ee5\u4e0b\u554f\u984c\uff1a \"),\n    @(\"\u75ab\u82d7\u63a5\u79cd\u8bc1\u4e66\u526f\u672c\", \"\u60a8\u7684\u75ab\u82d7\u63a5\u7a2e\u8b49\u660e\u526f\u672c\"),\n    @(\": \u6587\u4ef6\u4e0d\u6e05\u695a\", \": \u6587\u6a94\u4e0d\u6e05\u695a\"),\n    @(\"[\u6587\u4ef6 2]\", \"[\u6587\u6a94 2]\"),\n    @(\": [problem]\", \": [\u554f\u984c]\"),\n    @(\"\u8bf7\u5728 \", \"\u8acb\u5728 \"),\n    @(\" \u4e4b\u524d\u91cd\u65b0\u63d0\u4ea4\u4e0a\u8ff0\u6587\u4ef6\uff0c\u4ee5\u4fbf\u6211\u4eec\u8fdb\u884c\u5fc5\u8981\u7684\u5b89\u6392\u3002\", \" \u4e4b\u524d\u91cd\u65b0\u63d0\u4ea4\u4e0a\u8ff0\u6587\u6a94\uff0c\u4ee5\u4fbf\u6211\u5011\u9032\u884c\u5fc5\u8981\u7684\u5b89\u6392\u3002\"),\n    @(\"\u5982\u6709\u4efb\u4f55\u7591\u95ee\uff0c\u8bf7\u901a\u8fc7 \", \"\u5982\u6709\u4efb\u4f55\u7591\u554f\uff0c\u8acb\u901a\u904e \"),\n    @(\"[\u7535\u5b50\u90ae\u4ef6\u5730\u5740]\", \"[\u96fb\u5b50\u90f5\u4ef6\u5730\u5740]\"),\n    @(\"[WHATSAPP \u53f7\u7801]\", \"[WHATSAPP \u865f\u78bc]\"),\n    @(\" (WhatsApp) \u8054\u7cfb\u60a8\u7684\u533a\u57df\u7ecf\u7406 \", \" (WhatsApp) \u806f\u7e6b\u60a8\u7684\u5340\u57df\u7d93\u7406, \"),\n    @(\"[NAME]\", \"[\u59d3\u540d]\"),\n    @(\" \u3002 \", \"\u3002 \")\n)\n\nforeach ($pair in $replacements) {\n    $searchText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $null = $find.Execute($searchText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n}\n"}
